$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R1").Value = 5027
$ws.Range("S1").Value = 0.02472362494798537
$ws.Range("R2").Value = 4989
$ws.Range("S2").Value = 0.03436524775970404
$ws.Range("R3").Value = 5206
$ws.Range("S3").Value = 0.02890479927642864
$ws.Range("R4").Value = 4247
$ws.Range("S4").Value = 0.03042682410642743
$ws.Range("R5").Value = 5650
$ws.Range("S5").Value = 0.03003712341895088
$ws.Range("R6").Value = 5479
$ws.Range("S6").Value = 0.02699268327350968
$ws.Range("R7").Value = 5402
$ws.Range("S7").Value = 0.03082367591287649
$ws.Range("R8").Value = 6354
$ws.Range("S8").Value = 0.03088666625008012
$ws.Range("R9").Value = 5079
$ws.Range("S9").Value = 0.02841555478664241
$ws.Range("R10").Value = 5329
$ws.Range("S10").Value = 0.03237603515037423
$ws.Range("R11").Value = 5218
$ws.Range("S11").Value = 0.02555287587939577
$ws.Range("R12").Value = 5383
$ws.Range("S12").Value = 0.03704010792859991
$ws.Range("R13").Value = 5654
$ws.Range("S13").Value = 0.03681209818587049
$ws.Range("R14").Value = 4574
$ws.Range("S14").Value = 0.02411168792458445
